$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Bug fix: the "Pierre Papier Ciseau" quest's init code cell no longer ships
# a placeholder "Hello World" snippet.
$ws.Range("P2").ClearContents()

# Bug fix: correct the expected "Etat de l'eau" / "Moyenne" quest difficulty values.
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 2

# Scroll the view down a bit and move the selection to D3, like the author did.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D3").Select()
